# Auto-applied scheduled market-data refresh for Garuda_Profits workbook.
# Updates currentAveragePrice* / Leve* columns (H:N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 39265
$ws.Range("J63").Value = 39265
$ws.Range("L63").Value = 39265
$ws.Range("N63").Value = -40513
$ws.Range("H66").Value = 39265
$ws.Range("J66").Value = 39265
$ws.Range("L66").Value = 117795
$ws.Range("N66").Value = -124035
$ws.Range("H137").Value = 1356.6888
$ws.Range("I137").Value = 1106.0571
$ws.Range("J137").Value = 2233.9
$ws.Range("K137").Value = 3318.1713
$ws.Range("L137").Value = 6701.700000000001
$ws.Range("M137").Value = -768.1713
$ws.Range("N137").Value = -11801.7
$ws.Range("H138").Value = 1685.65
$ws.Range("I138").Value = 1313.7715
$ws.Range("J138").Value = 2206.28
$ws.Range("K138").Value = 3941.3145
$ws.Range("L138").Value = 6618.84
$ws.Range("M138").Value = 1198.6855
$ws.Range("N138").Value = -16898.84
$ws.Range("H141").Value = 1372.6582
$ws.Range("I141").Value = 748.4194
$ws.Range("J141").Value = 3649.2942
$ws.Range("K141").Value = 2245.2582
$ws.Range("L141").Value = 10947.8826
$ws.Range("M141").Value = 2934.7418
$ws.Range("N141").Value = -21307.8826

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 92
$ws.Range("I6").Value = 92
$ws.Range("K6").Value = 92
$ws.Range("M6").Value = 81
$ws.Range("H32").Value = 886.76
$ws.Range("I32").Value = 868.87366
$ws.Range("J32").Value = 1226.6
$ws.Range("K32").Value = 868.87366
$ws.Range("L32").Value = 1226.6
$ws.Range("M32").Value = -581.87366
$ws.Range("N32").Value = -1800.6
$ws.Range("H74").Value = 749.04877
$ws.Range("I74").Value = 624.62164
$ws.Range("K74").Value = 624.62164
$ws.Range("M74").Value = 249.37836
$ws.Range("H77").Value = 749.04877
$ws.Range("I77").Value = 624.62164
$ws.Range("K77").Value = 3123.1082
$ws.Range("M77").Value = 1244.8918
$ws.Range("H94").Value = 34915
$ws.Range("J94").Value = 34915
$ws.Range("L94").Value = 34915
$ws.Range("N94").Value = -36717
$ws.Range("H96").Value = 27585
$ws.Range("J96").Value = 27585
$ws.Range("L96").Value = 27585
$ws.Range("N96").Value = -33077

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2194.9048
$ws.Range("I20").Value = 1978.5862
$ws.Range("K20").Value = 1978.5862
$ws.Range("M20").Value = -1731.5862
$ws.Range("H57").Value = 51390
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 51390
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 51390
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -52830
$ws.Range("H135").Value = 25666.666
$ws.Range("J135").Value = 25666.666
$ws.Range("L135").Value = 25666.666
$ws.Range("N135").Value = -35806.666
$ws.Range("H136").Value = 51390
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 51390
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 51390
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -61590

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1062
$ws.Range("I16").Value = 885.7143
$ws.Range("J16").Value = 1267.6666
$ws.Range("K16").Value = 885.7143
$ws.Range("L16").Value = 1267.6666
$ws.Range("M16").Value = -598.7143
$ws.Range("N16").Value = -1841.6666
$ws.Range("H28").Value = 19319.6
$ws.Range("J28").Value = 19319.6
$ws.Range("L28").Value = 19319.6
$ws.Range("N28").Value = -19809.6
$ws.Range("H31").Value = 4388495
$ws.Range("I31").Value = 2021.64
$ws.Range("J31").Value = 12824020
$ws.Range("K31").Value = 2021.64
$ws.Range("L31").Value = 12824020
$ws.Range("M31").Value = -1726.64
$ws.Range("N31").Value = -12824610
$ws.Range("H34").Value = 4388495
$ws.Range("I34").Value = 2021.64
$ws.Range("J34").Value = 12824020
$ws.Range("K34").Value = 2021.64
$ws.Range("L34").Value = 12824020
$ws.Range("M34").Value = -1819.64
$ws.Range("N34").Value = -12824424
$ws.Range("H86").Value = 142858940
$ws.Range("I86").Value = 166668380
$ws.Range("J86").Value = 2300
$ws.Range("K86").Value = 166668380
$ws.Range("L86").Value = 2300
$ws.Range("M86").Value = -166667257
$ws.Range("N86").Value = -4546
$ws.Range("H88").Value = 24990
$ws.Range("J88").Value = 24990
$ws.Range("L88").Value = 24990
$ws.Range("N88").Value = -25802
$ws.Range("H89").Value = 142858940
$ws.Range("I89").Value = 166668380
$ws.Range("J89").Value = 2300
$ws.Range("K89").Value = 833341900
$ws.Range("L89").Value = 11500
$ws.Range("M89").Value = -833336284
$ws.Range("N89").Value = -22732
$ws.Range("H91").Value = 24990
$ws.Range("J91").Value = 24990
$ws.Range("L91").Value = 24990
$ws.Range("N91").Value = -27798
$ws.Range("H113").Value = 1062
$ws.Range("I113").Value = 885.7143
$ws.Range("J113").Value = 1267.6666
$ws.Range("K113").Value = 885.7143
$ws.Range("L113").Value = 1267.6666
$ws.Range("M113").Value = 1284.2857
$ws.Range("N113").Value = -5607.6666
$ws.Range("H132").Value = 1861.5111
$ws.Range("I132").Value = 1898.7646
$ws.Range("J132").Value = 1746.3636
$ws.Range("K132").Value = 5696.293799999999
$ws.Range("L132").Value = 5239.0908
$ws.Range("M132").Value = -3166.293799999999
$ws.Range("N132").Value = -10299.0908

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1536985.9
$ws.Range("J131").Value = 2025311.4
$ws.Range("L131").Value = 6075934.199999999
$ws.Range("N131").Value = -6086014.199999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9276658
$ws.Range("I70").Value = 11864428
$ws.Range("J70").Value = 3813.6667
$ws.Range("K70").Value = 11864428
$ws.Range("L70").Value = 3813.6667
$ws.Range("M70").Value = -11864158
$ws.Range("N70").Value = -4353.6667
$ws.Range("H73").Value = 9276658
$ws.Range("I73").Value = 11864428
$ws.Range("J73").Value = 3813.6667
$ws.Range("K73").Value = 11864428
$ws.Range("L73").Value = 3813.6667
$ws.Range("M73").Value = -11863492
$ws.Range("N73").Value = -5685.6667

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1104.0834
$ws.Range("I40").Value = 1057.1428
$ws.Range("J40").Value = 1169.8
$ws.Range("K40").Value = 1057.1428
$ws.Range("L40").Value = 1169.8
$ws.Range("M40").Value = -921.1428000000001
$ws.Range("N40").Value = -1441.8
$ws.Range("H104").Value = 29979.5
$ws.Range("J104").Value = 29979.5
$ws.Range("L104").Value = 29979.5
$ws.Range("N104").Value = -36967.5
$ws.Range("H132").Value = 7222.6763
$ws.Range("I132").Value = 10527.714
$ws.Range("J132").Value = 1883.7693
$ws.Range("K132").Value = 31583.142
$ws.Range("L132").Value = 5651.3079
$ws.Range("M132").Value = -29053.142
$ws.Range("N132").Value = -10711.3079

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 31502.143
$ws.Range("J125").Value = 31502.143
$ws.Range("L125").Value = 31502.143
$ws.Range("N125").Value = -41342.143
$ws.Range("H132").Value = 899.65
$ws.Range("I132").Value = 877.94116
$ws.Range("J132").Value = 1022.6667
$ws.Range("K132").Value = 2633.82348
$ws.Range("L132").Value = 3068.0001
$ws.Range("M132").Value = -103.82348
$ws.Range("N132").Value = -8128.0001
